# Add MTZ (Miller-Tucker-Zemlin) constraints to avoid loops.
#
# The sIntersections sheet lists one row per MTZ ordering-variable
# (intersection, v01VisitIntersection) and the sPaths sheet lists one row
# per directed path constraint (pOriginIntersection, pDestinationIntersection,
# v01TravelPath). Both tables are re-numbered/re-indexed and grow with extra
# rows so every intersection/arc used by the MTZ subtour-elimination
# constraints is represented.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # sIntersections
$ws2 = $wb.Worksheets.Item(2)   # sPaths

# ---------------------------------------------------------------------------
# sIntersections: column A (intersection id) for rows 2-16; column B
# (v01VisitIntersection) stays 1 for every row. Rows 2-10 already exist and
# get re-numbered; rows 11-16 are brand new.
# ---------------------------------------------------------------------------
$intersectionIds = @(1, 2, 5, 8, 9, 14, 19, 20, 24, 25, 27, 29, 30, 31, 38)

$row = 2
foreach ($id in $intersectionIds) {
    $ws1.Cells.Item($row, 1).Value = $id
    $ws1.Cells.Item($row, 2).Value = 1
    $row = $row + 1
}

# The source table carries empty (text-typed) placeholder cells across
# columns C:I for every data row. Recreate that on the newly added rows
# 11-16: "'" forces an empty-text cell instead of clearing it outright, then
# copy the same blank style used elsewhere in the table onto those cells so
# the new rows keep the exact same shape/formatting as the rest of the table.
$ws1.Range("C11:I16").Value = "'"
$blankStyle = $ws1.Cells.Item(2, 3).Style
$ws1.Range("C11:I16").Style = $blankStyle

# ---------------------------------------------------------------------------
# sPaths: columns A (pOriginIntersection) / B (pDestinationIntersection) for
# rows 2-15; column C (v01TravelPath) stays 1 for every row. Rows 2-9 already
# exist and get re-pointed to the renumbered intersections; rows 10-15 are
# brand new arcs needed by the added MTZ constraints.
# ---------------------------------------------------------------------------
$paths = @(
    @(2, 5), @(5, 14), @(14, 31), @(24, 20), @(29, 27), @(25, 38), @(8, 1), @(1, 2),
    @(9, 8), @(19, 9), @(20, 19), @(27, 25), @(30, 29), @(31, 30)
)

$row = 2
foreach ($p in $paths) {
    $ws2.Cells.Item($row, 1).Value = $p[0]
    $ws2.Cells.Item($row, 2).Value = $p[1]
    $ws2.Cells.Item($row, 3).Value = 1
    $row = $row + 1
}

Write-Output "Added MTZ constraint rows: sIntersections now $($intersectionIds.Count + 1) rows, sPaths now $($paths.Count + 1) rows"
